$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 37.98512966666667
$ws.Cells.Item(2, 8).Value = 113.955389
$ws.Cells.Item(2, 9).Value = 0.5085441461893128
$ws.Cells.Item(2, 10).Value = 0.5085441461893129
$ws.Cells.Item(2, 13).Value = 0.2901893333333334
$ws.Cells.Item(2, 14).Value = 0.870568
$ws.Cells.Item(2, 15).Value = 0.03429389578125064
$ws.Cells.Item(2, 16).Value = 0.03429389578125064
$ws.Cells.Item(2, 17).Value = 11.02287945455022
$ws.Cells.Item(2, 18).Value = 99.205915090952
$ws.Cells.Item(2, 19).Value = 0.01743995994958138
$ws.Cells.Item(2, 20).Value = 0.01743995994958138

$ws.Cells.Item(3, 7).Value = 37.98512966666667
$ws.Cells.Item(3, 8).Value = 113.955389
$ws.Cells.Item(3, 9).Value = 0.5085441461893128
$ws.Cells.Item(3, 10).Value = 0.5085441461893129
$ws.Cells.Item(3, 15).Value = 0.8402845891331153
$ws.Cells.Item(3, 16).Value = 0.8402845891331153
$ws.Cells.Item(3, 17).Value = 270.0875920488031
$ws.Cells.Item(3, 18).Value = 2430.788328439228
$ws.Cells.Item(3, 19).Value = 0.4273218089367376
$ws.Cells.Item(3, 20).Value = 0.4273218089367377

$ws.Cells.Item(4, 7).Value = 37.98512966666667
$ws.Cells.Item(4, 8).Value = 113.955389
$ws.Cells.Item(4, 9).Value = 0.5085441461893128
$ws.Cells.Item(4, 10).Value = 0.5085441461893129
$ws.Cells.Item(4, 15).Value = 0.1254215150856341
$ws.Cells.Item(4, 16).Value = 0.1254215150856341
$ws.Cells.Item(4, 17).Value = 40.31347883642456
$ws.Cells.Item(4, 18).Value = 362.821309527821
$ws.Cells.Item(4, 19).Value = 0.06378237730299381
$ws.Cells.Item(4, 20).Value = 0.06378237730299381

$ws.Cells.Item(5, 9).Value = 0.1771904651558058
$ws.Cells.Item(5, 10).Value = 0.1771904651558058
$ws.Cells.Item(5, 13).Value = 0.2901893333333334
$ws.Cells.Item(5, 14).Value = 0.870568
$ws.Cells.Item(5, 15).Value = 0.03429389578125064
$ws.Cells.Item(5, 16).Value = 0.03429389578125064
$ws.Cells.Item(5, 17).Value = 3.840667821159111
$ws.Cells.Item(5, 18).Value = 34.566010390432
$ws.Cells.Item(5, 19).Value = 0.006076551345484526
$ws.Cells.Item(5, 20).Value = 0.006076551345484527

$ws.Cells.Item(6, 9).Value = 0.1771904651558058
$ws.Cells.Item(6, 10).Value = 0.1771904651558058
$ws.Cells.Item(6, 15).Value = 0.8402845891331153
$ws.Cells.Item(6, 16).Value = 0.8402845891331153
$ws.Cells.Item(6, 19).Value = 0.1488904172117518
$ws.Cells.Item(6, 20).Value = 0.1488904172117519

$ws.Cells.Item(7, 9).Value = 0.1771904651558058
$ws.Cells.Item(7, 10).Value = 0.1771904651558058
$ws.Cells.Item(7, 15).Value = 0.1254215150856341
$ws.Cells.Item(7, 16).Value = 0.1254215150856341
$ws.Cells.Item(7, 19).Value = 0.02222349659856942
$ws.Cells.Item(7, 20).Value = 0.02222349659856942

$ws.Cells.Item(8, 8).Value = 70.421093
$ws.Cells.Item(8, 9).Value = 0.3142653886548814
$ws.Cells.Item(8, 10).Value = 0.3142653886548814
$ws.Cells.Item(8, 13).Value = 0.2901893333333334
$ws.Cells.Item(8, 14).Value = 0.870568
$ws.Cells.Item(8, 15).Value = 0.03429389578125064
$ws.Cells.Item(8, 16).Value = 0.03429389578125064
$ws.Cells.Item(8, 17).Value = 6.811816676758222
$ws.Cells.Item(8, 18).Value = 61.306350090824
$ws.Cells.Item(8, 19).Value = 0.01077738448618473
$ws.Cells.Item(8, 20).Value = 0.01077738448618473

$ws.Cells.Item(9, 8).Value = 70.421093
$ws.Cells.Item(9, 9).Value = 0.3142653886548814
$ws.Cells.Item(9, 10).Value = 0.3142653886548814
$ws.Cells.Item(9, 15).Value = 0.8402845891331153
$ws.Cells.Item(9, 16).Value = 0.8402845891331153
$ws.Cells.Item(9, 19).Value = 0.2640723629846258
$ws.Cells.Item(9, 20).Value = 0.2640723629846258

$ws.Cells.Item(10, 8).Value = 70.421093
$ws.Cells.Item(10, 9).Value = 0.3142653886548814
$ws.Cells.Item(10, 10).Value = 0.3142653886548814
$ws.Cells.Item(10, 15).Value = 0.1254215150856341
$ws.Cells.Item(10, 16).Value = 0.1254215150856341
$ws.Cells.Item(10, 17).Value = 24.91254926340856
$ws.Cells.Item(10, 19).Value = 0.03941564118407087
$ws.Cells.Item(10, 20).Value = 0.03941564118407087
